# Auto-generated edit script applying the diff to Carbuncle_Profits (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 276
$ws.Range("I2").Value = 211.42857
$ws.Range("J2").Value = 426.66666
$ws.Range("K2").Value = 211.42857
$ws.Range("L2").Value = 426.66666
$ws.Range("M2").Value = -98.42857000000001
$ws.Range("N2").Value = -652.66666
# Row 9
$ws.Range("H9").Value = 198.22223
$ws.Range("I9").Value = 197.71428
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 197.71428
$ws.Range("L9").Value = 200
$ws.Range("M9").Value = -28.71428
$ws.Range("N9").Value = -538
# Row 98
$ws.Range("H98").Value = 1116.875
$ws.Range("I98").Value = 1116.875
$ws.Range("K98").Value = 1116.875
$ws.Range("M98").Value = 381.125
# Row 112
$ws.Range("H112").Value = 1025.6818
$ws.Range("J112").Value = 1024.4736
$ws.Range("L112").Value = 3073.4208
$ws.Range("N112").Value = -5289.4208
# Row 122
$ws.Range("H122").Value = 1116.875
$ws.Range("I122").Value = 1116.875
$ws.Range("K122").Value = 3350.625
$ws.Range("M122").Value = -900.625
# Row 132
$ws.Range("H132").Value = 44092.832
$ws.Range("I132").Value = 45870.824
$ws.Range("K132").Value = 137612.472
$ws.Range("M132").Value = -135082.472
# Row 138
$ws.Range("H138").Value = 12348575
$ws.Range("I138").Value = 1333.475
$ws.Range("J138").Value = 24394664
$ws.Range("K138").Value = 4000.425
$ws.Range("L138").Value = 73183992
$ws.Range("M138").Value = 1139.575
$ws.Range("N138").Value = -73194272

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1507.5555
$ws.Range("I2").Value = 1285.7273
$ws.Range("J2").Value = 1856.1428
$ws.Range("K2").Value = 1285.7273
$ws.Range("L2").Value = 1856.1428
$ws.Range("M2").Value = -1172.7273
$ws.Range("N2").Value = -2082.1428
# Row 32
$ws.Range("H32").Value = 3963
$ws.Range("I32").Value = 2640.1516
$ws.Range("J32").Value = 13663.889
$ws.Range("K32").Value = 2640.1516
$ws.Range("L32").Value = 13663.889
$ws.Range("M32").Value = -2353.1516
$ws.Range("N32").Value = -14237.889
# Row 45
$ws.Range("H45").Value = 2435.5386
$ws.Range("I45").Value = 1073.5555
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 1073.5555
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -696.5554999999999
$ws.Range("N45").Value = -6254
# Row 102
$ws.Range("H102").Value = 2400
$ws.Range("I102").Value = 2400
$ws.Range("K102").Value = 2400
$ws.Range("M102").Value = -778
# Row 110
$ws.Range("H110").Value = 100974.336
$ws.Range("I110").Value = 150458.25
$ws.Range("J110").Value = 2006.5
$ws.Range("K110").Value = 150458.25
$ws.Range("L110").Value = 2006.5
$ws.Range("M110").Value = -148413.25
$ws.Range("N110").Value = -6096.5
# Row 116
$ws.Range("H116").Value = 1507.5555
$ws.Range("I116").Value = 1285.7273
$ws.Range("J116").Value = 1856.1428
$ws.Range("K116").Value = 1285.7273
$ws.Range("L116").Value = 1856.1428
$ws.Range("M116").Value = 1008.2727
$ws.Range("N116").Value = -6444.1428

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1507.5555
$ws.Range("I3").Value = 1285.7273
$ws.Range("J3").Value = 1856.1428
$ws.Range("K3").Value = 1285.7273
$ws.Range("L3").Value = 1856.1428
$ws.Range("M3").Value = -1171.7273
$ws.Range("N3").Value = -2084.1428
# Row 99
$ws.Range("H99").Value = 1471.0526
$ws.Range("I99").Value = 883.3333
$ws.Range("K99").Value = 883.3333
$ws.Range("M99").Value = 614.6667
# Row 107
$ws.Range("H107").Value = 984.7143
$ws.Range("I107").Value = 998.8333
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 998.8333
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 921.1667
$ws.Range("N107").Value = -4740

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 309
$ws.Range("I22").Value = 327.625
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 327.625
$ws.Range("L22").Value = 160
$ws.Range("M22").Value = 22.375
$ws.Range("N22").Value = -860
# Row 31
$ws.Range("H31").Value = 3685.8875
$ws.Range("I31").Value = 976.7143
$ws.Range("J31").Value = 5144.673
$ws.Range("K31").Value = 976.7143
$ws.Range("L31").Value = 5144.673
$ws.Range("M31").Value = -681.7143
$ws.Range("N31").Value = -5734.673
# Row 34
$ws.Range("H34").Value = 3685.8875
$ws.Range("I34").Value = 976.7143
$ws.Range("J34").Value = 5144.673
$ws.Range("K34").Value = 976.7143
$ws.Range("L34").Value = 5144.673
$ws.Range("M34").Value = -774.7143
$ws.Range("N34").Value = -5548.673
# Row 134
$ws.Range("H134").Value = 4123.7144
$ws.Range("I134").Value = 4594.3335
$ws.Range("K134").Value = 13783.0005
$ws.Range("M134").Value = -11248.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 966.45917
$ws.Range("I68").Value = 677.9138
$ws.Range("J68").Value = 1384.85
$ws.Range("K68").Value = 2033.7414
$ws.Range("L68").Value = 4154.549999999999
$ws.Range("M68").Value = -1222.7414
$ws.Range("N68").Value = -5776.549999999999
# Row 71
$ws.Range("H71").Value = 966.45917
$ws.Range("I71").Value = 677.9138
$ws.Range("J71").Value = 1384.85
$ws.Range("K71").Value = 6101.224200000001
$ws.Range("L71").Value = 12463.65
$ws.Range("M71").Value = -2045.224200000001
$ws.Range("N71").Value = -20575.65
# Row 137
$ws.Range("H137").Value = 2330.647
$ws.Range("I137").Value = 1797.3158
$ws.Range("J137").Value = 3006.2
$ws.Range("K137").Value = 5391.9474
$ws.Range("L137").Value = 9018.599999999999
$ws.Range("M137").Value = -291.9474
$ws.Range("N137").Value = -19218.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1349.6875
$ws.Range("I2").Value = 1697.8334
$ws.Range("J2").Value = 1140.8
$ws.Range("K2").Value = 1697.8334
$ws.Range("L2").Value = 1140.8
$ws.Range("M2").Value = -1584.8334
$ws.Range("N2").Value = -1366.8
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
# Row 122
$ws.Range("H122").Value = 94841.21000000001
$ws.Range("I122").Value = 126806.25
$ws.Range("J122").Value = 14928.625
$ws.Range("K122").Value = 380418.75
$ws.Range("L122").Value = 44785.875
$ws.Range("M122").Value = -377968.75
$ws.Range("N122").Value = -49685.875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1461.3636
$ws.Range("I113").Value = 1527.5
$ws.Range("K113").Value = 4582.5
$ws.Range("M113").Value = -2412.5

Write-Host "Edits applied"